$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("⬛", "📘", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("🟥", "📕", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("🟧", "📙", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("🟩", "📗", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("noir", "bleu", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
